$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 57: 7 May 2020 (Excel serial date 43958) Covid-19 Turkey data
$ws.Range("A57").Value = 43958
$ws.Range("B57").Value = 30395
$ws.Range("C57").Value = 1977
$ws.Range("D57").Value = 57
$ws.Range("E57").Value = 4782

# Grow the worksheet table ("Table3") to cover the newly added row
$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E57"))

# Update the window scroll position / selection to match the saved view state
$win = $wb.Windows.Item(1)
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("D56").Select() | Out-Null
